# Auto-generated Excel COM-interop script that applies the "Updated cryptos list" edit.
# Updates Price (D) and Volume(1h) (E) columns for rows 2-51 of the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.790.71"
$ws.Range("E2").Value = "  -0.96%  "
$ws.Range("D3").Value = "2.309.79"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'303.07"
$ws.Range("E5").Value = "  -1.81%  "
$ws.Range("D6").Value = "'99.77"
$ws.Range("E6").Value = "  -4.08%  "
$ws.Range("E7").Value = "  -3.61%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.504"
$ws.Range("E9").Value = "  -2.91%  "
$ws.Range("D10").Value = "'34.82"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").Value = "'0.0793"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  -3.30%  "
$ws.Range("D14").Value = "2.667.27"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "'15.69"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").Value = "2.301.30"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'0.806"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "42.697.61"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "'11.56"
$ws.Range("E20").Value = "  -3.41%  "
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").Value = "'67.98"
$ws.Range("D23").Value = "'235.22"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  -2.20%  "
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "'24.98"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("D29").Value = "'34.69"
$ws.Range("E29").Value = "  -4.30%  "
$ws.Range("D30").Value = "'164.11"
$ws.Range("E30").Value = "  +1.61%  "
$ws.Range("E31").Value = "  -4.20%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  +0.01%  "
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  -4.30%  "
$ws.Range("E34").Value = "  -4.84%  "
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("D36").Value = "'16.85"
$ws.Range("E36").Value = "  -7.44%  "
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  -5.75%  "
$ws.Range("E41").Value = "  -3.54%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "1.968.94"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  -2.78%  "
$ws.Range("D45").Value = "'18.57"
$ws.Range("E45").Value = "  -1.25%  "
$ws.Range("D46").Value = "'10.22"
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("D47").Value = "'2.88"
$ws.Range("E47").Value = "  -6.13%  "
$ws.Range("D48").Value = "'55.53"
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("D49").Value = "2.534.01"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -3.15%  "
$ws.Range("E51").Value = "  +0.76%  "
